$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G6").Value = 2120
$ws.Range("D8").Value = 1270
$ws.Range("G8").Value = 750
$ws.Range("C9").Value = 1400
$ws.Range("D9").Value = 1250

$ws.Range("C14").Select()
